$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the autofilter and its backing defined name (_xlnm._FilterDatabase) ---
$ws.AutoFilterMode = $false
foreach ($n in @($wb.Names)) {
    $n.Delete()
}

# --- Drop the now-unused trailing columns (F:J) and empty trailing rows (2:9) ---
$ws.Range("F1:J1").EntireColumn.Delete()
$ws.Range("A2:A9").EntireRow.Delete()

# --- Add the new "siglas" column header in E1 ---
$ws.Range("E1").Value = "siglas"

# --- Re-style the sheet: base font becomes Arial 10 instead of Calibri 11 ---
$ws.Range("A1:E1").Font.Name = "Arial"

# Header cells A1:D1 -> Arial 12, bold removed, black, centered
$hdr = $ws.Range("A1:D1")
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 12
$hdr.Font.Bold = $false
$hdr.Font.Color = 0
$hdr.HorizontalAlignment = -4108  # xlCenter

# E1 -> Arial 12, no fill (previous yellow highlight removed)
$e1 = $ws.Range("E1")
$e1.Font.Name = "Arial"
$e1.Font.Size = 12
$e1.Interior.Pattern = -4142 # xlNone

# Remove the yellow fill still present on the header row defaults (B1/D1 originally)
$hdr.Interior.Pattern = -4142 # xlNone

# --- Column widths (characters) ---
$ws.Columns.Item(1).ColumnWidth = 6.333333
$ws.Columns.Item(2).ColumnWidth = 19.666667
$ws.Columns.Item(3).ColumnWidth = 19.5
$ws.Columns.Item(4).ColumnWidth = 22.0
$ws.Columns.Item(5).ColumnWidth = 6.5

# --- Row height for header row ---
$ws.Rows.Item(1).RowHeight = 15

# --- View: select whole sheet, land on E5 ---
$ws.Range("A1:XFD1048576").Select()
$excel.ActiveWindow.RangeSelection.Application.Goto($ws.Range("E5"))
$ws.Range("A1:XFD1048576").Select()

# --- Page setup ---
$ps = $ws.PageSetup
$ps.LeftMargin = $excel.InchesToPoints(0.7875)
$ps.RightMargin = $excel.InchesToPoints(0.7875)
$ps.TopMargin = $excel.InchesToPoints(1.05277777777778)
$ps.BottomMargin = $excel.InchesToPoints(1.05277777777778)
$ps.HeaderMargin = $excel.InchesToPoints(0.7875)
$ps.FooterMargin = $excel.InchesToPoints(0.7875)
$ps.FirstPageNumber = 1
$ps.CenterHeader = "&""Times New Roman,Normal""&12&A"
$ps.CenterFooter = "&""Times New Roman,Normal""&12Página &P"

Write-Host "done"
